$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 0.9211691253470349
$ws.Range("C2").Value = 0.1834521536900695
$ws.Range("D2").Value = 0.6008320174162094
$ws.Range("E2").Value = 0.2299847524138805
$ws.Range("G2").Value = 0.649463987586941
$ws.Range("H2").Value = 0.762768681843454
$ws.Range("J2").Value = 0.1102173371868922
$ws.Range("M2").Value = 0.4223793129907349
$ws.Range("O2").Value = 2.815072590898609
$ws.Range("B3").Value = 0.8204206955624045
$ws.Range("C3").Value = 0.1601982326561142
$ws.Range("D3").Value = 0.5976860526324259
$ws.Range("E3").Value = 0.2301952801200073
$ws.Range("G3").Value = 0.6535320393179305
$ws.Range("H3").Value = 0.76977997774312
$ws.Range("J3").Value = 0.1111422963292164
$ws.Range("M3").Value = 0.3942834481389568
$ws.Range("O3").Value = 2.838037381637051
$ws.Range("B4").Value = 0.7584867392029082
$ws.Range("C4").Value = 0.1458617823548138
$ws.Range("D4").Value = 0.5960564809138162
$ws.Range("E4").Value = 0.2304423709180412
$ws.Range("G4").Value = 0.6566080749850585
$ws.Range("H4").Value = 0.7745265418449421
$ws.Range("J4").Value = 0.1117662674607942
$ws.Range("M4").Value = 0.3771157369693086
$ws.Range("O4").Value = 2.854276389589984
$ws.Range("B5").Value = 0.7332310023204798
$ws.Range("C5").Value = 0.1400052480499028
$ws.Range("D5").Value = 0.595468440687128
$ws.Range("E5").Value = 0.2305726989930932
$ws.Range("G5").Value = 0.6580066410284218
$ws.Range("H5").Value = 0.7765717961737479
$ws.Range("J5").Value = 0.1120346285358416
$ws.Range("M5").Value = 0.3701410795864959
$ws.Range("O5").Value = 2.861430817576675
$ws.Range("B6").Value = 0.7290363142411707
$ws.Range("C6").Value = 0.1390319217937304
$ws.Range("D6").Value = 0.5953753906817099
$ws.Range("E6").Value = 0.2305961300953321
$ws.Range("G6").Value = 0.6582476227742404
$ws.Range("H6").Value = 0.776918111126605
$ws.Range("J6").Value = 0.1120800404229403
$ws.Range("M6").Value = 0.3689842424496774
$ws.Range("O6").Value = 2.862651207389561
$ws.Range("B7").Value = 0.7581461989780109
$ws.Range("C7").Value = 0.1457828565309853
$ws.Range("D7").Value = 0.596048242491733
$ws.Range("E7").Value = 0.2304440085595161
$ws.Range("G7").Value = 0.6566263496678673
$ws.Range("H7").Value = 0.7745536755347615
$ws.Range("J7").Value = 0.1117698296373106
$ws.Range("M7").Value = 0.3770215874431813
$ws.Range("O7").Value = 2.854370703956334
$ws.Range("B8").Value = 0.8864472714478779
$ws.Range("C8").Value = 0.1754465543072854
$ws.Range("D8").Value = 0.5996846331648413
$ws.Range("E8").Value = 0.2300328919794694
$ws.Range("G8").Value = 0.6507464304273611
$ws.Range("H8").Value = 0.7650944905681456
$ws.Range("J8").Value = 0.110524628713037
$ws.Range("M8").Value = 0.412674799500671
$ws.Range("O8").Value = 2.822546484051117
$ws.Range("B9").Value = 1.137406482943163
$ws.Range("C9").Value = 0.2331392715973664
$ws.Range("D9").Value = 0.6092105083856723
$ws.Range("E9").Value = 0.2301614935988177
$ws.Range("G9").Value = 0.6438201194639248
$ws.Range("H9").Value = 0.7500515237307752
$ws.Range("J9").Value = 0.1085276939229374
$ws.Range("M9").Value = 0.4832380893547708
$ws.Range("O9").Value = 2.777147571141455
$ws.Range("B10").Value = 1.321346367459455
$ws.Range("C10").Value = 0.2752204067390664
$ws.Range("D10").Value = 0.6176684292192647
$ws.Range("E10").Value = 0.2308261254953479
$ws.Range("G10").Value = 0.6415613816203773
$ws.Range("H10").Value = 0.7411409060671303
$ws.Range("J10").Value = 0.1073321190103798
$ws.Range("M10").Value = 0.5354634860852769
$ws.Range("O10").Value = 2.754219333576515
$ws.Range("B11").Value = 1.404920561335075
$ws.Range("C11").Value = 0.294295077892059
$ws.Range("D11").Value = 0.6218329816871631
$ws.Range("E11").Value = 0.2312523737492143
$ws.Range("G11").Value = 0.6411533592697225
$ws.Range("H11").Value = 0.7375530833143529
$ws.Range("J11").Value = 0.1068472709433621
$ws.Range("M11").Value = 0.5593030588673855
$ws.Range("O11").Value = 2.746065207111712
$ws.Range("B12").Value = 1.436552239962168
$ws.Range("C12").Value = 0.301508029824106
$ws.Range("D12").Value = 0.6234555373125374
$ws.Range("E12").Value = 0.2314316004633525
$ws.Range("G12").Value = 0.6410883243661658
$ws.Range("H12").Value = 0.7362615006559139
$ws.Range("J12").Value = 0.1066721662507391
$ws.Range("M12").Value = 0.5683419725318402
$ws.Range("O12").Value = 2.743305728816694
$ws.Range("B13").Value = 1.429740530300933
$ws.Range("C13").Value = 0.2999550520403602
$ws.Range("D13").Value = 0.6231040668336902
$ws.Range("E13").Value = 0.2313922084458504
$ws.Range("G13").Value = 0.6410983451668102
$ws.Range("H13").Value = 0.7365366827391711
$ws.Range("J13").Value = 0.1067095001395479
$ws.Range("M13").Value = 0.5663947807893521
$ws.Range("O13").Value = 2.743885413845135
$ws.Range("B14").Value = 1.407523249064354
$ws.Range("C14").Value = 0.294888698715738
$ws.Range("D14").Value = 0.6219655581369921
$ws.Range("E14").Value = 0.2312667617927602
$ws.Range("G14").Value = 0.6411462133674348
$ws.Range("H14").Value = 0.7374454797963921
$ws.Range("J14").Value = 0.1068326946355711
$ws.Range("M14").Value = 0.5600464698636358
$ws.Range("O14").Value = 2.745831597469618
$ws.Range("B15").Value = 1.393912386736417
$ws.Range("C15").Value = 0.2917840693709479
$ws.Range("D15").Value = 0.6212741164239901
$ws.Range("E15").Value = 0.2311922421900761
$ws.Range("G15").Value = 0.6411871978435073
$ws.Range("H15").Value = 0.7380108784208375
$ws.Range("J15").Value = 0.1069092616183411
$ws.Range("M15").Value = 0.5561594168474215
$ws.Range("O15").Value = 2.747066477326428
$ws.Range("B16").Value = 1.315882433860395
$ws.Range("C16").Value = 0.2739724246588651
$ws.Range("D16").Value = 0.6174026401438937
$ws.Range("E16").Value = 0.2308007623040176
$ws.Range("G16").Value = 0.6416005461309879
$ws.Range("H16").Value = 0.7413847564081948
$ws.Range("J16").Value = 0.1073649934755281
$ws.Range("M16").Value = 0.533907128968778
$ws.Range("O16").Value = 2.754798112009411
$ws.Range("B17").Value = 1.267986624987714
$ws.Range("C17").Value = 0.2630277928330713
$ws.Range("D17").Value = 0.6151087677531848
$ws.Range("E17").Value = 0.2305923357503232
$ws.Range("G17").Value = 0.6420130676653173
$ws.Range("H17").Value = 0.7435738498747071
$ws.Range("J17").Value = 0.1076596943193344
$ws.Range("M17").Value = 0.520276779362348
$ws.Range("O17").Value = 2.760124895703626
$ws.Range("B18").Value = 1.240428808280967
$ws.Range("C18").Value = 0.2567263297308955
$ws.Range("D18").Value = 0.6138192339810757
$ws.Range("E18").Value = 0.2304841168111516
$ws.Range("G18").Value = 0.6423086249839685
$ws.Range("H18").Value = 0.7448767836685448
$ws.Range("J18").Value = 0.107834753613858
$ws.Range("M18").Value = 0.5124447028111234
$ws.Range("O18").Value = 2.763402892003143
$ws.Range("B19").Value = 1.231096635136225
$ws.Range("C19").Value = 0.2545916755343001
$ws.Range("D19").Value = 0.6133877468815285
$ws.Range("E19").Value = 0.2304494788303018
$ws.Range("G19").Value = 0.6424186942801953
$ws.Range("H19").Value = 0.7453254590453611
$ws.Range("J19").Value = 0.1078949795139827
$ws.Range("M19").Value = 0.5097942379760099
$ws.Range("O19").Value = 2.764549518216455
$ws.Range("B20").Value = 1.273086201105286
$ws.Range("C20").Value = 0.2641935328869351
$ws.Range("D20").Value = 0.6153498663551318
$ws.Range("E20").Value = 0.2306133160986086
$ws.Range("G20").Value = 0.6419631188276185
$ws.Range("H20").Value = 0.7433362807740167
$ws.Range("J20").Value = 0.1076277479074541
$ws.Range("M20").Value = 0.5217269552175026
$ws.Range("O20").Value = 2.759535677162376
$ws.Range("B21").Value = 1.414049450987079
$ws.Range("C21").Value = 0.2963770898008988
$ws.Range("D21").Value = 0.6222987305011145
$ws.Range("E21").Value = 0.2313031250126016
$ws.Range("G21").Value = 0.6411297218885181
$ws.Range("H21").Value = 0.7371767237204239
$ws.Range("J21").Value = 0.1067962787541354
$ws.Range("M21").Value = 0.5619108168841507
$ws.Range("O21").Value = 2.745251037743003
$ws.Range("B22").Value = 1.506082803223308
$ws.Range("C22").Value = 0.3173511995671845
$ws.Range("D22").Value = 0.6271055734740969
$ws.Range("E22").Value = 0.231857798228571
$ws.Range("G22").Value = 0.6411067546730465
$ws.Range("H22").Value = 0.7335419440111934
$ws.Range("J22").Value = 0.1063023953885178
$ws.Range("M22").Value = 0.5882394683120538
$ws.Range("O22").Value = 2.737829327481876
$ws.Range("B23").Value = 1.456972000438554
$ws.Range("C23").Value = 0.3061625157930052
$ws.Range("D23").Value = 0.6245158085109495
$ws.Range("E23").Value = 0.2315522574962046
$ws.Range("G23").Value = 0.6410711472637303
$ws.Range("H23").Value = 0.7354461032892914
$ws.Range("J23").Value = 0.1065614552535266
$ws.Range("M23").Value = 0.5741814489049517
$ws.Range("O23").Value = 2.741614954886643
$ws.Range("B24").Value = 1.27078075077759
$ws.Range("C24").Value = 0.2636665305878978
$ws.Range("D24").Value = 0.6152407745823893
$ws.Range("E24").Value = 0.2306037947245514
$ws.Range("G24").Value = 0.641985518833863
$ws.Range("H24").Value = 0.7434435474806094
$ws.Range("J24").Value = 0.1076421733379185
$ws.Range("M24").Value = 0.5210713176238926
$ws.Range("O24").Value = 2.759801391483791
$ws.Range("B25").Value = 1.069589028631469
$ws.Range("C25").Value = 0.2175845197946842
$ws.Range("D25").Value = 0.6063771925217907
$ws.Range("E25").Value = 0.2300265992551331
$ws.Range("G25").Value = 0.6451985580098949
$ws.Range("H25").Value = 0.7537452596826313
$ws.Range("J25").Value = 0.109020252159997
$ws.Range("M25").Value = 0.4640808068768223
$ws.Range("O25").Value = 2.787602357991346
